$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.189922480620155
$ws.Range("C2").Value = 0.5697674418604651
$ws.Range("J2").Value = 0.01550387596899225
$ws.Range("P2").Value = 0.1317829457364341
$ws.Range("S2").Value = 0.09302325581395349
$ws.Range("C3").Value = 0.0261437908496732
$ws.Range("J3").Value = 0.07843137254901961
$ws.Range("P3").Value = 0.7124183006535948
$ws.Range("S3").Value = 0.1830065359477124
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.6904761904761905
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.08653846153846154
$ws.Range("D6").Value = 0.004807692307692308
$ws.Range("F6").Value = 0.125
$ws.Range("J6").Value = 0.2307692307692308
$ws.Range("O6").Value = 0.009615384615384616
$ws.Range("Q6").Value = 0.2067307692307692
$ws.Range("R6").Value = 0.05288461538461538
$ws.Range("S6").Value = 0.2836538461538461
$ws.Range("B7").Value = 0.08641975308641975
$ws.Range("D7").Value = 0.0308641975308642
$ws.Range("F7").Value = 0.04938271604938271
$ws.Range("J7").Value = 0.1358024691358025
$ws.Range("O7").Value = 0.03703703703703703
$ws.Range("Q7").Value = 0.2098765432098765
$ws.Range("R7").Value = 0.1049382716049383
$ws.Range("S7").Value = 0.345679012345679
$ws.Range("B8").Value = 0.125
$ws.Range("D8").Value = 0.02163461538461538
$ws.Range("F8").Value = 0.06490384615384616
$ws.Range("J8").Value = 0.08653846153846154
$ws.Range("O8").Value = 0.02163461538461538
$ws.Range("Q8").Value = 0.1778846153846154
$ws.Range("R8").Value = 0.09615384615384616
$ws.Range("S8").Value = 0.40625
$ws.Range("B9").Value = 0.07894736842105263
$ws.Range("D9").Value = 0.01973684210526316
$ws.Range("F9").Value = 0.07894736842105263
$ws.Range("J9").Value = 0.1842105263157895
$ws.Range("O9").Value = 0.01973684210526316
$ws.Range("Q9").Value = 0.1776315789473684
$ws.Range("R9").Value = 0.09210526315789473
$ws.Range("S9").Value = 0.3486842105263158
$ws.Range("B10").Value = 0.1067518248175182
$ws.Range("D10").Value = 0.02281021897810219
$ws.Range("F10").Value = 0.06478102189781022
$ws.Range("J10").Value = 0.1049270072992701
$ws.Range("O10").Value = 0.02463503649635037
$ws.Range("Q10").Value = 0.2363138686131387
$ws.Range("R10").Value = 0.0948905109489051
$ws.Range("S10").Value = 0.3448905109489051
$ws.Range("F11").Value = 0.004629629629629629
$ws.Range("G11").Value = 0.1481481481481481
$ws.Range("J11").Value = 0.09259259259259259
$ws.Range("K11").Value = 0.1805555555555556
$ws.Range("L11").Value = 0.5740740740740741
$ws.Range("G12").Value = 0.8015873015873016
$ws.Range("J12").Value = 0.1587301587301587
$ws.Range("K12").Value = 0.007936507936507936
$ws.Range("L12").Value = 0.007936507936507936
$ws.Range("S12").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01382488479262673
$ws.Range("H15").Value = 0.1382488479262673
$ws.Range("I15").Value = 0.05529953917050692
$ws.Range("J15").Value = 0.3456221198156682
$ws.Range("K15").Value = 0.08294930875576037
$ws.Range("O15").Value = 0.1105990783410138
$ws.Range("S15").Value = 0.2534562211981567
$ws.Range("F16").Value = 0.02994011976047904
$ws.Range("H16").Value = 0.1796407185628743
$ws.Range("I16").Value = 0.0658682634730539
$ws.Range("J16").Value = 0.4431137724550898
$ws.Range("K16").Value = 0.08383233532934131
$ws.Range("M16").Value = 0.01197604790419162
$ws.Range("O16").Value = 0.0718562874251497
$ws.Range("S16").Value = 0.1137724550898204
$ws.Range("F17").Value = 0.01376146788990826
$ws.Range("H17").Value = 0.2064220183486239
$ws.Range("I17").Value = 0.07339449541284404
$ws.Range("J17").Value = 0.4334862385321101
$ws.Range("K17").Value = 0.06422018348623854
$ws.Range("M17").Value = 0.03669724770642202
$ws.Range("N17").Value = 0.002293577981651376
$ws.Range("O17").Value = 0.05275229357798165
$ws.Range("S17").Value = 0.1169724770642202
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.2096774193548387
$ws.Range("I18").Value = 0.1021505376344086
$ws.Range("J18").Value = 0.4193548387096774
$ws.Range("K18").Value = 0.07526881720430108
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("N18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.04838709677419355
$ws.Range("S18").Value = 0.1021505376344086
$ws.Range("F19").Value = 0.0224390243902439
$ws.Range("H19").Value = 0.2282926829268293
$ws.Range("I19").Value = 0.07707317073170732
$ws.Range("J19").Value = 0.3707317073170732
$ws.Range("K19").Value = 0.09463414634146342
$ws.Range("M19").Value = 0.02341463414634146
$ws.Range("N19").Value = 0.001951219512195122
$ws.Range("O19").Value = 0.07609756097560975
$ws.Range("S19").Value = 0.1053658536585366
